$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.22"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "22.10"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.358"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05853"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.395"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.378"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8138"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.012"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1423"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.04115"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCXBestin24h"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07341"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03032"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.175"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09393"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001592"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04818"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0005891"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006090"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.004082"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0009862"
$ws.Range("E21").Value = "20BitKanKAN"
$ws.Range("B22").Value = "NitroEx"
$ws.Range("C22").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0001410"
$ws.Range("E22").Value = "21NitroExNTX"
$ws.Range("B23").Value = "LEO"
$ws.Range("C23").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.707"
$ws.Range("E23").Value = "22LEOLEO"
$ws.Range("B24").Value = "BTSEToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.231"
$ws.Range("E24").Value = "23BTSETokenBTSE"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3248"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002472"
$ws.Range("E27").Value = "26UpBotsUBXTWorstin24h"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006411"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003001"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.005071"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005643"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7221"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.08856"
$ws.Range("E48").Value = "47BOLOBOLO"
